$d = $word.ActiveDocument

# The paragraph that currently ends with "...under en threshhold. " is the
# last paragraph in the document body (before the sectPr). We will append
# three new paragraphs after it:
#   1. A new "Overskrift2" (Heading 2) paragraph: "15-01-2021"
#   2. A new bulleted "Listeafsnit" paragraph about Friday
#   3. A new bulleted "Listeafsnit" paragraph about the droplet/moon pixel art,
#      which also receives the relocated "_GoBack" bookmark.

$pMainscreen = $d.Paragraphs.Last

# --- 1. "15-01-2021" heading ---------------------------------------------
$endOfMainscreen = $pMainscreen.Range.Duplicate
$endOfMainscreen.Collapse(0)
$endOfMainscreen.InsertParagraphAfter()
$pHeading = $d.Paragraphs.Last
$pHeading.Style = "Overskrift2"
$pHeading.Range.Text = "15-01-2021"

# --- 2. Friday bullet paragraph -------------------------------------------
$endOfHeading = $pHeading.Range.Duplicate
$endOfHeading.Collapse(0)
$endOfHeading.InsertParagraphAfter()
$pFredag = $d.Paragraphs.Last
$pFredag.Style = "Listeafsnit"
$pFredag.Range.Text = "Fredag er en kort dag, og jeg hjalp mine klassekammerater en del med dels seriel forbindelse og noget timing problemer."

# Give this paragraph a fresh bulleted list (matches the "Symbol" bullet
# list already used elsewhere in the document).
$bulletGallery = $word.ListGalleries.Item(1)
$bulletTemplate = $bulletGallery.ListTemplates.Item(1)
$pFredag.Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate)

# --- 3. Droplet / moon bullet paragraph ------------------------------------
$endOfFredag = $pFredag.Range.Duplicate
$endOfFredag.Collapse(0)
$endOfFredag.InsertParagraphAfter()
$pDraabe = $d.Paragraphs.Last
# New paragraph inherits the "Listeafsnit" style + the same numbered list
# (numId) from the preceding paragraph automatically.

# Write the full sentence first so later bookmark math is not at the very
# end of the document content (collapsed ranges right at document end are
# unreliable for Bookmarks.Add in this runtime).
$pDraabe.Range.Text = " Jeg har også lavet nogle kompliceret grafiske figurer så som  en vand dråbe, sky og måne."

# Relocate the existing "_GoBack" bookmark so that it now sits right after
# "...dråbe, sky" and right before " og måne." (re-adding a bookmark with
# the same name moves it rather than duplicating it).
$findRng = $d.Content.Duplicate
$findRng.Find.Execute("dråbe, sky", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkPoint = $d.Range($findRng.End, $findRng.End)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint)
